# GroupVol.xlsx edit:
#  - Remove the "Table" column (A) from the generated INSERT formulas in column F
#  - Replace the VolID GUID values in column C with a fresh set of GUIDs
#  - Normalise the style of column C back to the default (style index 0)
#  - Move the active selection to C15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New VolID values (column C) per row, and the formulas (column F) that no
# longer reference the "Table" column (A).
$rows = @(
    @{ Row = 2;  VolId = "dba53101-f9b2-4dc0-85e7-11d472fd99cd" },
    @{ Row = 3;  VolId = "589178b4-aa4c-4276-a516-9460fa7714d3" },
    @{ Row = 4;  VolId = "293fe520-7e35-444a-8955-f02a911fed1c" },
    @{ Row = 5;  VolId = "34fb4310-9790-4b80-84cc-8c899f0308f7" },
    @{ Row = 6;  VolId = "6bd2cd4d-e58c-4c9f-8195-8c4f0ae9af1d" },
    @{ Row = 7;  VolId = "63aaaec0-0222-4cc6-b748-668f8675fb40" },
    @{ Row = 8;  VolId = "55f42ce3-4b5d-44ad-bb75-a0719ac5bab3" },
    @{ Row = 9;  VolId = "dba53101-f9b2-4dc0-85e7-11d472fd99cd" },
    @{ Row = 10; VolId = "46e0eab8-9d77-4a4d-a642-bed325a80ba2" },
    @{ Row = 11; VolId = "ec311095-16c4-4ea1-a9bc-9ddcda3b9b62" }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Update the VolID text in column C.
    $ws.Cells.Item($r, 3).Value = $item.VolId

    # Column C previously carried an explicit (bold-ish) style; restore it to
    # the sheet's plain default style (same as columns D/E) by copying that
    # formatting across instead of assigning a brand-new style record.
    $ws.Cells.Item($r, 4).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)

    # Rebuild the INSERT statement formula in column F without the
    # "[" &A$1 &"]," piece and without the A<row> value argument.
    $ws.Cells.Item($r, 6).Formula = '="INSERT INTO "&A$2&" (["' + '&B$1&"],["&C$1&"],["&D$1&"],["&E$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'' ,''" & D' + $r + ' & "'',''" & E' + $r + ' & "'')"'
}

$excel.CutCopyMode = 0

# Move the active selection.
$ws.Range("C15").Select()
